# "videos script creation BDD"
# Adds a second worksheet ("Feuil2") holding two small reference tables
# (realisateurs / films) backed by Excel ListObjects, mirroring the
# structure described by the target diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes right after the existing "Feuil1" tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuil2"

# ---- Tableau1: realisateurs (D8:E11) ----
$ws2.Range("D8").Value = "realisateur_id"
$ws2.Range("E8").Value = "realisateur_nom"
$ws2.Range("D9").Value = 1
$ws2.Range("E9").Value = "Toto"
$ws2.Range("D10").Value = 2
$ws2.Range("E10").Value = "Tata"
$ws2.Range("D11").Value = 3
$ws2.Range("E11").Value = "Titi"

# ---- Tableau2: films (G15:I19) ----
$ws2.Range("G15").Value = "film_id"
$ws2.Range("H15").Value = "film_titre"
$ws2.Range("I15").Value = "realisateur_id"
$ws2.Range("G16").Value = 1
$ws2.Range("H16").Value = "Léon"
$ws2.Range("I16").Value = 1
$ws2.Range("G17").Value = 2
$ws2.Range("H17").Value = "E.T"
$ws2.Range("I17").Value = 3
$ws2.Range("G18").Value = 3
$ws2.Range("H18").Value = "ça"
$ws2.Range("I18").Value = 2
$ws2.Range("G19").Value = 4
$ws2.Range("H19").Value = "Identity"
$ws2.Range("I19").Value = 1

# Column widths (best-effort match of the author's manual resize).
$ws2.Columns.Item(4).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(5).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(7).ColumnWidth = 8.666666666666666
$ws2.Columns.Item(8).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(9).ColumnWidth = 14.666666666666666

# Turn both ranges into real Excel Tables.
$ws2.ListObjects.Add(1, $ws2.Range("D8:E11"), $null, 1) | Out-Null
$ws2.ListObjects.Add(1, $ws2.Range("G15:I19"), $null, 1) | Out-Null

# Rename highest index first -- renaming low-to-high clobbers later
# tables back to their auto-generated default name in this host.
$ws2.ListObjects.Item(2).Name = "Tableau2"
$ws2.ListObjects.Item(1).Name = "Tableau1"

# Match the author's final view state: zoomed in on Feuil2, looking at
# the films table.
$ws2.Activate() | Out-Null
$ws2.Range("E15").Select() | Out-Null
$excel.ActiveWindow.Zoom = 145

Write-Host "Feuil2 created with Tableau1/Tableau2"
